$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Tabulate (indent) the pseudo code lines, in top-to-bottom order so the
# shared-strings table gets the new entries appended in that same order.

# Row 36: "playerVal + trainVal <= MAX_VAL" gets a 4-space indent
$ws.Range("B36").Value = "    playerVal + trainVal <= MAX_VAL"

# Row 37: "then playerVal = playerVal + trainVal" gets a 4-space indent
$ws.Range("B37").Value = "    then playerVal = playerVal + trainVal"

# Row 39: "begin if" gets a 4-space indent
$ws.Range("B39").Value = "    begin if"

# Row 40: "playerVal + trainVal > MAX_VAL" gets an 8-space indent
$ws.Range("B40").Value = "        playerVal + trainVal > MAX_VAL"

# Row 41: "then overByMaxVal = (playerVal + trainVal) - maxVal" gets an 8-space indent;
# the note in D41 moves over to E41 (unchanged text)
$ws.Range("B41").Value = "        then overByMaxVal = (playerVal + trainVal) - maxVal"
$noteText = $ws.Range("D41").Value2
$ws.Range("D41").ClearContents()
$ws.Range("E41").Value = $noteText

# Row 42: "end if" gets a 3-space indent
$ws.Range("B42").Value = "   end if"

# Row 43: "playerVal = (playerVal + trainVal) - overByVal" gets a 3-space indent
$ws.Range("B43").Value = "   playerVal = (playerVal + trainVal) - overByVal"

# Reflect the resulting selection position
$ws.Range("B43").Select()
